$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.802.42'
$ws.Range("E2").Value = '  +1.05%  '
$ws.Range("D3").Value = '1.649.67'
$ws.Range("E3").Value = '  +1.35%  '
$ws.Range("E4").Value = '  +0.64%  '
$ws.Range("D5").Value = '216.81'
$ws.Range("E6").Value = '  +0.91%  '
$ws.Range("E7").Value = '  +0.74%  '
$ws.Range("D8").Value = '0.252'
$ws.Range("E8").Value = '  +1.45%  '
$ws.Range("E9").Value = '  +0.79%  '
$ws.Range("E10").Value = '  +2.07%  '
$ws.Range("D11").Value = '0.0844'
$ws.Range("E11").Value = '  +0.14%  '
$ws.Range("E12").Value = '  +1.39%  '
$ws.Range("D13").Value = '1.656.95'
$ws.Range("E13").Value = '  +1.66%  '
$ws.Range("D14").Value = '4.19'
$ws.Range("E14").Value = '  +1.36%  '
$ws.Range("D16").Value = '65.48'
$ws.Range("E16").Value = '  +0.77%  '
$ws.Range("D17").Value = '26.817.93'
$ws.Range("E17").Value = '  +1.04%  '
$ws.Range("E18").Value = '  +0.72%  '
$ws.Range("D19").Value = '218.45'
$ws.Range("E19").Value = '  +1.82%  '
$ws.Range("E20").Value = '  +0.62%  '
$ws.Range("E21").Value = '  +1.68%  '
$ws.Range("D22").Value = '2.41'
$ws.Range("E22").Value = '  +15.88%  '
$ws.Range("E23").Value = '  +0.28%  '
$ws.Range("D24").Value = '9.50'
$ws.Range("E24").Value = '  +2.32%  '
$ws.Range("D25").Value = '146.65'
$ws.Range("E25").Value = '  -1.20%  '
$ws.Range("E26").Value = '  +0.50%  '
$ws.Range("E27").Value = '  -0.06%  '
$ws.Range("E28").Value = '  +3.89%  '
$ws.Range("E29").Value = '  +1.26%  '
$ws.Range("E30").Value = '  +1.82%  '
$ws.Range("E31").Value = '  +2.12%  '
$ws.Range("E32").Value = '  +0.30%  '
$ws.Range("E33").Value = '  +1.75%  '
$ws.Range("D34").Value = '1.286.57'
$ws.Range("E34").Value = '  +3.92%  '
$ws.Range("E35").Value = '  +3.09%  '
$ws.Range("E36").Value = '  +3.14%  '
$ws.Range("D37").Value = '0.0178'
$ws.Range("E37").Value = '  +2.43%  '
$ws.Range("E38").Value = '  +5.93%  '
$ws.Range("D39").Value = '0.828'
$ws.Range("E39").Value = '  +4.15%  '
$ws.Range("E40").Value = '  +0.78%  '
$ws.Range("D41").Value = '0.815'
$ws.Range("E41").Value = '  +2.06%  '
$ws.Range("E42").Value = '  -0.81%  '
$ws.Range("E43").Value = '  +2.57%  '
$ws.Range("D44").Value = '1.789.63'
$ws.Range("E44").Value = '  +1.41%  '
$ws.Range("D45").Value = '92.09'
$ws.Range("E45").Value = '  -0.91%  '
$ws.Range("D46").Value = '59.72'
$ws.Range("E46").Value = '  +8.75%  '
$ws.Range("E47").Value = '  +1.47%  '
$ws.Range("E48").Value = '  +1.29%  '
$ws.Range("D49").Value = '7.74'
$ws.Range("E49").Value = '  +3.19%  '
$ws.Range("D50").Value = '0.0971'
$ws.Range("E50").Value = '  +1.98%  '
$ws.Range("E51").Value = '  +0.64%  '
